$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump: 5.0.0 -> 6.0.0
$ws.Range("B3").Value = "6.0.0"

# Date refresh
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank; now filled in
$ws.Range("B9").Value = "Alvearie Team"

# Row 10 was "Contact" / "No display for ContactDetail"; becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Row 11 duplicated the old "Contact" row - remove it entirely, shifting the remaining rows up
$ws.Rows.Item(11).Delete()
